$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-9) got reshuffled: some rows were replaced with
# the content of other rows (a cyclic rotation among rows 2,3,4,5,6, and a
# swap between rows 8 and 9). Only columns A, B, D, E, F, G, H, Q, R change;
# all other columns are identical across the affected rows so they need no
# changes.

# New values for each affected row, taken from the diff (after state).
$rowData = @{
    2 = @{ A = 97650292; B = 5135;  D = "LC"; E = 105930; F = "Vågbandad barkbock"; G = "Semanotus undatus";     H = "(Linnaeus, 1758)";   Q = 403705.050704394;  R = 6794737.908215457 }
    3 = @{ A = 97650293; B = 77506; D = "NT"; E = 6425;   F = "Garnlav";             G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 403710.6419448711; R = 6794695.894706693 }
    4 = @{ A = 97650291; B = 5113;  D = "LC"; E = 100526; F = "Bronshjon";           G = "Callidium coriaceum";   H = "Paykull, 1800";      Q = 403714.9324539425; R = 6794746.778207967 }
    5 = @{ A = 97650299; B = 77506; D = "NT"; E = 6425;   F = "Garnlav";             G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 403880.3826524244; R = 6794050.283030285 }
    6 = @{ A = 97650294; B = 90676; D = "NT"; E = 5966;   F = "Motaggsvamp";         G = "Sarcodon squamosus";    H = "(Schaeff.) Quél.";   Q = 403863.9880530759; R = 6794102.706117956 }
    8 = @{ A = 97650298; B = 96334; D = "VU"; E = 220787; F = "Knärot";              G = "Goodyera repens";       H = "(L.) R. Br.";        Q = 403840.5463236904; R = 6794038.864283022 }
    9 = @{ A = 97650301; B = 90653; D = "LC"; E = 4364;   F = "Dropptaggsvamp";      G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."; Q = 403960.8920370748; R = 6793787.235077787 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $vals.A   # A
    $ws.Cells.Item($r, 2).Value = $vals.B   # B
    $ws.Cells.Item($r, 4).Value = $vals.D   # D
    $ws.Cells.Item($r, 5).Value = $vals.E   # E
    $ws.Cells.Item($r, 6).Value = $vals.F   # F
    $ws.Cells.Item($r, 7).Value = $vals.G   # G
    $ws.Cells.Item($r, 8).Value = $vals.H   # H
    $ws.Cells.Item($r, 17).Value = $vals.Q  # Q
    $ws.Cells.Item($r, 18).Value = $vals.R  # R
}
